$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number, new Price (D) text, new Volume(1h) (E) text,
# and whether the Price text must be forced to stay text (it parses as a
# plain number like "1.00" / "527.36", which Excel would otherwise
# auto-convert to a numeric value and mangle, e.g. dropping trailing zeros).
$updates = @(
    [PSCustomObject]@{ Row = 2; Price = "60.766.69"; Volume = "  -1.53%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 3; Price = "2.905.00"; Volume = "  -2.47%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 4; Price = "1.00"; Volume = "  -0.02%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 5; Price = "527.36"; Volume = "  -2.43%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 6; Price = "143.99"; Volume = "  -4.95%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 7; Price = "0.999"; Volume = "  +0.01%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 8; Price = "0.548"; Volume = "  -3.37%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 9; Price = "2.912.69"; Volume = "  -2.61%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 10; Price = "0.108"; Volume = "  -4.42%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 11; Price = "6.02"; Volume = "  -2.22%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 12; Price = "0.359"; Volume = "  -2.67%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 13; Price = "3.413.80"; Volume = "  -2.55%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 14; Price = "0.129"; Volume = "  +3.26%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 15; Price = "60.702.79"; Volume = "  -1.71%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 16; Price = "22.61"; Volume = "  -5.52%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 17; Price = "2.908.11"; Volume = "  -2.41%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 18; Price = "0.0000140"; Volume = "  -3.95%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 19; Price = "4.96"; Volume = "  -3.89%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 20; Price = "11.59"; Volume = "  -3.77%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 21; Price = "351.17"; Volume = "  -7.94%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 22; Price = "6.50"; Volume = "  -3.13%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 23; Price = "1.00"; Volume = "  -0.05%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 24; Price = "5.73"; Volume = "  +1.29%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 25; Price = "64.96"; Volume = "  -1.44%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 26; Price = "0.451"; Volume = "  -4.21%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 27; Price = "0.178"; Volume = "  -5.96%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 28; Price = "0.998"; Volume = "  -1.48%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 29; Price = "7.85"; Volume = "  -3.70%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 30; Price = "0.0₃0860"; Volume = "  -8.28%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 31; Price = "1.00"; Volume = "  +0.02%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 32; Price = "1.68"; Volume = "  -2.47%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 33; Price = "19.56"; Volume = "  -4.49%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 34; Price = "152.08"; Volume = "  -5.01%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 35; Price = "4.34"; Volume = "  -4.89%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 36; Price = "5.57"; Volume = "  -6.00%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 37; Price = "0.993"; Volume = "  -7.02%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 38; Price = "1.20"; Volume = "  -5.69%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 39; Price = "37.60"; Volume = "  +0.04%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 40; Price = "1.47"; Volume = "  -4.88%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 41; Price = "3.72"; Volume = "  -4.49%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 42; Price = "2.292.98"; Volume = "  -5.49%  "; ForceText = $false }
    [PSCustomObject]@{ Row = 43; Price = "0.651"; Volume = "  -3.08%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 44; Price = "0.0581"; Volume = "  -1.52%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 45; Price = "20.38"; Volume = "  -7.32%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 46; Price = "0.997"; Volume = "  +0.01%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 47; Price = "4.96"; Volume = "  -4.26%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 48; Price = "0.0238"; Volume = "  -2.46%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 49; Price = "10.34"; Volume = "  -0.84%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 50; Price = "0.0917"; Volume = "  -3.66%  "; ForceText = $true }
    [PSCustomObject]@{ Row = 51; Price = "18.40"; Volume = "  -7.06%  "; ForceText = $true }
)

foreach ($u in $updates) {
    $priceCell = $ws.Cells.Item($u.Row, 4)   # column D = Price
    $volumeCell = $ws.Cells.Item($u.Row, 5)  # column E = Volume(1h)

    if ($u.ForceText) {
        # Prefix with an apostrophe so Excel stores the literal text instead of
        # re-parsing it as a number, then reset the style back to Normal so no
        # stray number-format / quote-prefix style sticks to the cell.
        $priceCell.Value = "'" + $u.Price
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $u.Price
    }

    $volumeCell.Value = $u.Volume
}
